# Add a new BOM row (row 11) for the 2x10 shrouded header connector.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Shrouded Header"
$ws.Range("F11").Value = "CNN HEADDER 2.54mm 10POS GOLD"
$ws.Range("G11").Value = "Sullins Connector Solutions"
$ws.Range("H11").Value = "SBH11-PBPC-D05-ST-BK"
$ws.Range("I11").Value = "S9169-ND"

$ws.Range("K11").Value = 0.64
$ws.Range("L11").Value = 0.445
$ws.Range("M11").Value = 0.2985
$ws.Range("N11").Value = 0.2159

# I11 picks up the "Digikey part number" look (Arial 9) used by the rest
# of column I, matching I10's formatting.
$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The cursor ends up on B11 after typing in the new row.
$ws.Range("B11").Select()

$wb.Save()
